# Update work-hour tracking workbook (commit: "update 31-10, check old github page for more info")
#
# Summary of the edit being reproduced:
#  - "Teamlid 1" logged an extra 30 minutes on the css3 activity (row 32):
#      D32: 60 -> 90   (this cascades through the SUM/formulas on every sheet)
#  - The description of that same css3 activity (C32) gained more detail.
#  - The active sheet/tab moved from "Logboek Totalen" to "Teamlid 1", and the
#    selected cell on each sheet changed to reflect where the author was working.

$wb = $excel.ActiveWorkbook

$wsLogboek  = $wb.Worksheets.Item("Logboek Totalen")
$wsTeam1    = $wb.Worksheets.Item("Teamlid 1")
$wsTeam2    = $wb.Worksheets.Item("Teamlid 2")

# Expand the description of the css3 styling activity.
$wsTeam1.Range("C32").Value = "leren css3 en toevoegen speciale styling voor form posts met css3 (grotere eerste letter post, afwisselende background-color posts, animatie bij hover posts en text-shadow)"

# Extra 30 minutes logged for that activity; D43/tot_teamlid1 and the Logboek
# Totalen formulas recompute automatically from this single edit.
$wsTeam1.Range("D32").Value = 90

# Leave "Logboek Totalen" selection untouched (still D7), just move the active
# tab to "Teamlid 1" and update the selections that moved as the author
# continued working: Teamlid 1 -> C30, Teamlid 2 -> C15.
# ("Teamlid 1" is selected last so it ends up the active/selected tab.)
[void]$wsTeam2.Range("C15").Select()
[void]$wsTeam1.Range("C30").Select()
